$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.849.72'
$ws.Range("E2").Value = '  -4.27%  '
$ws.Range("D3").Value = '3.016.26'
$ws.Range("E3").Value = '  -4.54%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.21'
$ws.Range("E5").Value = '  -7.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.02'
$ws.Range("E6").Value = '  -10.63%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.011.60'
$ws.Range("E8").Value = '  -4.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("E10").Value = '  -4.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.95'
$ws.Range("E11").Value = '  -12.34%  '
$ws.Range("E12").Value = '  -5.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000218'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.16'
$ws.Range("E14").Value = '  -10.34%  '
$ws.Range("D15").Value = '3.475.90'
$ws.Range("E15").Value = '  -4.95%  '
$ws.Range("D16").Value = '61.922.48'
$ws.Range("E16").Value = '  -4.15%  '
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").Value = '3.023.19'
$ws.Range("E18").Value = '  -4.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.41'
$ws.Range("E19").Value = '  -6.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.30'
$ws.Range("E20").Value = '  -9.64%  '
$ws.Range("E21").Value = '  -8.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.680'
$ws.Range("E22").Value = '  -5.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.89'
$ws.Range("E23").Value = '  -8.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.46'
$ws.Range("E24").Value = '  -2.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.61'
$ws.Range("E25").Value = '  -9.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.62'
$ws.Range("E27").Value = '  -7.91%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.93'
$ws.Range("E28").Value = '  -10.36%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.25'
$ws.Range("E30").Value = '  -5.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.80'
$ws.Range("E31").Value = '  -16.44%  '
$ws.Range("E32").Value = '  -5.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("E33").Value = '  -11.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '55.53'
$ws.Range("E34").Value = '  +2.88%  '
$ws.Range("E35").Value = '  -5.90%  '
$ws.Range("E36").Value = '  -6.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '454.83'
$ws.Range("E37").Value = '  -18.06%  '
$ws.Range("D38").Value = '3.036.98'
$ws.Range("E38").Value = '  -4.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0383'
$ws.Range("E39").Value = '  -12.33%  '
$ws.Range("E40").Value = '  -7.52%  '
$ws.Range("E41").Value = '  -9.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.84'
$ws.Range("E42").Value = '  -5.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.47'
$ws.Range("E43").Value = '  -10.67%  '
$ws.Range("E45").Value = '  -9.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.95'
$ws.Range("E46").Value = '  -12.53%  '
$ws.Range("D47").Value = '0.0₃0507'
$ws.Range("E47").Value = '  -2.62%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.105'
$ws.Range("E48").Value = '  -3.02%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.57'
$ws.Range("E49").Value = '  -7.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '115.14'
$ws.Range("E50").Value = '  -4.88%  '
$ws.Range("E51").Value = '  -9.48%  '
